$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9: date shifts from 2023-08-11 (45149) to 2023-08-14 (45152)
$ws.Range("A9").Value = 45152

# Fill in attendance for row 9 (Present/Absent per person)
$ws.Range("B9").Value = "Present"
$ws.Range("C9").Value = "Present"
$ws.Range("D9").Value = "Present"
$ws.Range("E9").Value = "Absent"
$ws.Range("F9").Value = "Absent"
$ws.Range("G9").Value = "Absent"
$ws.Range("H9").Value = "Present"
$ws.Range("I9").Value = "Absent"

# Update row 10: date shifts from 2023-08-12 (45150) to 2023-08-15 (45153)
$ws.Range("A10").Value = 45153

# Update the selected cell in the sheet view
$ws.Range("F10").Select()

# Adjust column widths: column G should match column H's existing width
# (~12.7109375 chars). Column H is left untouched.
$ws.Range("G1").ColumnWidth = 11.8
